# Daily attendance processing - 2025-11-01 18:49:06
# Normalize the "Recorded By" (column G) values so that "System" (when
# paired with other recorders) is moved from the front of the list to the
# end, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)   # strip leading "System, "
        $cell.Value2 = $rest + ", System"
    }
}
